$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.496.33"
$ws.Range("E2").Value = "  -1.55%  "

$ws.Range("D3").Value = "1.851.69"
$ws.Range("E3").Value = "  -0.75%  "

$ws.Range("D4").Value = "'0.9992"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "'243.35"
$ws.Range("E5").Value = "  -1.16%  "

$ws.Range("D6").Value = "'0.6430"
$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").Value = "'48.27"
$ws.Range("E8").Value = "  +2.20%  "

$ws.Range("D9").Value = "'0.2989"
$ws.Range("E9").Value = "  -0.82%  "

$ws.Range("D10").Value = "'0.07462"
$ws.Range("E10").Value = "  -0.39%  "

$ws.Range("D11").Value = "'24.26"
$ws.Range("E11").Value = "  -0.83%  "

$ws.Range("D12").Value = "'0.07628"
$ws.Range("E12").Value = "  -0.73%  "

$ws.Range("D13").Value = "1.852.68"
$ws.Range("E13").Value = "  -0.86%  "

$ws.Range("D14").Value = "'5.019"
$ws.Range("E14").Value = "  -1.09%  "

$ws.Range("D15").Value = "'0.6839"
$ws.Range("E15").Value = "  -1.03%  "

$ws.Range("D16").Value = "'83.60"
$ws.Range("E16").Value = "  -0.76%  "

$ws.Range("D17").Value = "'0.000009524"
$ws.Range("E17").Value = "  +0.62%  "

$ws.Range("D18").Value = "'6.142"
$ws.Range("E18").Value = "  +0.50%  "

$ws.Range("D19").Value = "29.512.71"
$ws.Range("E19").Value = "  -1.61%  "

$ws.Range("D20").Value = "2.070.74"
$ws.Range("E20").Value = "  -2.73%  "

$ws.Range("D21").Value = "'235.33"
$ws.Range("E21").Value = "  -2.40%  "

$ws.Range("D22").Value = "'12.55"
$ws.Range("E22").Value = "  -1.29%  "

$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("D24").Value = "'7.678"
$ws.Range("E24").Value = "  +2.79%  "

$ws.Range("E25").Value = "  -0.15%  "

$ws.Range("D26").Value = "'157.13"
$ws.Range("E26").Value = "  -1.60%  "

$ws.Range("D27").Value = "'0.1412"
$ws.Range("E27").Value = "  -1.05%  "

$ws.Range("D28").Value = "'8.476"
$ws.Range("E28").Value = "  -1.50%  "

$ws.Range("D29").Value = "'17.75"
$ws.Range("E29").Value = "  -1.68%  "

$ws.Range("E30").Value = "  -1.41%  "

$ws.Range("D31").Value = "'0.06000"
$ws.Range("E31").Value = "  -2.77%  "

$ws.Range("E32").Value = "  -1.64%  "

$ws.Range("D33").Value = "'4.123"
$ws.Range("E33").Value = "  -1.42%  "

$ws.Range("D34").Value = "'4.070"
$ws.Range("E34").Value = "  -1.53%  "

$ws.Range("D35").Value = "'1.865"
$ws.Range("E35").Value = "  -0.62%  "

$ws.Range("D36").Value = "'1.175"
$ws.Range("E36").Value = "  +0.70%  "

$ws.Range("D37").Value = "'0.7187"
$ws.Range("E37").Value = "  -2.50%  "

$ws.Range("E38").Value = "  -0.24%  "

$ws.Range("D39").Value = "'2.799"
$ws.Range("E39").Value = "  -2.55%  "

$ws.Range("D40").Value = "'0.01776"
$ws.Range("E40").Value = "  -1.78%  "

$ws.Range("D41").Value = "1.198.63"
$ws.Range("E41").Value = "  -2.06%  "

$ws.Range("D42").Value = "'0.9081"
$ws.Range("E42").Value = "  -2.44%  "

$ws.Range("D43").Value = "'6.173"
$ws.Range("E43").Value = "  -1.94%  "

$ws.Range("E44").Value = "  -0.15%  "

$ws.Range("D45").Value = "2.007.79"
$ws.Range("E45").Value = "  -1.59%  "

$ws.Range("D46").Value = "'101.86"
$ws.Range("E46").Value = "  -0.51%  "

$ws.Range("D47").Value = "'66.33"
$ws.Range("E47").Value = "  -0.59%  "

$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.315"
$ws.Range("E48").Value = "  +8.52%  "

$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.00000000120"
$ws.Range("E49").Value = "  -2.49%  "

$ws.Range("D50").Value = "'0.4032"
$ws.Range("E50").Value = "  -1.78%  "

$ws.Range("D51").Value = "'9.079"
$ws.Range("E51").Value = "  -3.13%  "
